$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "SP.URB.TOTL:GIB"
$ws.Range("A2").Value = "SP.POP.TOTL:GIB:cor-value"
$ws.Range("B2").Value = 1
$ws.Range("A3").Value = "SP.POP.TOTL:GIB:p-value"
$ws.Range("B3").Value = 0

$headerRange = $ws.Range("B1:B1,A2:A2,A3:A3")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2
$headerRange.Borders.ColorIndex = -4105
